$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure column D (Price) cells stay text even if the new value looks numeric
# (mirrors the original workbook where these are stored as inline strings).

$ws.Range('D2').NumberFormat = "@"
$ws.Range('D2').Value = '67.972.17'
$ws.Range('E2').Value = '  +1.17%  '
$ws.Range('D3').NumberFormat = "@"
$ws.Range('D3').Value = '2.619.95'
$ws.Range('E3').Value = '  +0.40%  '
$ws.Range('D4').NumberFormat = "@"
$ws.Range('D4').Value = '1.00'
$ws.Range('E4').Value = '  -0.20%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '596.02'
$ws.Range('E5').Value = '  +0.43%  '
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D6').Value = '153.37'
$ws.Range('E6').Value = '  +0.40%  '
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D7').Value = '1.00'
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.544'
$ws.Range('E8').Value = '  -1.56%  '
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D9').Value = '2.622.38'
$ws.Range('E9').Value = '  +0.41%  '
$ws.Range('E10').Value = '  +9.92%  '
$ws.Range('E11').Value = '  -0.48%  '
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D12').Value = '5.22'
$ws.Range('E12').Value = '  +0.63%  '
$ws.Range('D13').NumberFormat = "@"
$ws.Range('D13').Value = '0.347'
$ws.Range('E13').Value = '  +0.02%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '27.57'
$ws.Range('E14').Value = '  -0.26%  '
$ws.Range('E15').Value = '  +4.39%  '
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D16').Value = '3.094.79'
$ws.Range('E16').Value = '  +0.23%  '
$ws.Range('D17').NumberFormat = "@"
$ws.Range('D17').Value = '67.826.91'
$ws.Range('E17').Value = '  +1.34%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '2.616.76'
$ws.Range('E18').Value = '  +0.01%  '
$ws.Range('E19').Value = '  +2.72%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '366.90'
$ws.Range('E20').Value = '  -0.02%  '
$ws.Range('E21').Value = '  +0.36%  '
$ws.Range('D22').NumberFormat = "@"
$ws.Range('D22').Value = '4.22'
$ws.Range('E22').Value = '  -1.58%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '4.80'
$ws.Range('E23').Value = '  -0.96%  '
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D24').Value = '2.07'
$ws.Range('E24').Value = '  +0.75%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '72.76'
$ws.Range('E25').Value = '  +8.92%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '9.94'
$ws.Range('E27').Value = '  -1.65%  '
$ws.Range('B28').Value = 'WrappedeETH'
$ws.Range('C28').Value = 'https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth'
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D28').Value = '2.750.65'
$ws.Range('E28').Value = '  +0.41%  '
$ws.Range('B29').Value = 'PEPE'
$ws.Range('C29').Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D29').Value = '0.0000104'
$ws.Range('E29').Value = '  +2.87%  '
$ws.Range('E30').Value = '  -0.30%  '
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D31').Value = '569.78'
$ws.Range('E31').Value = '  -2.73%  '
$ws.Range('B32').Value = 'InternetComputer(DFINITY)'
$ws.Range('C32').Value = 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.85'
$ws.Range('E32').Value = '  +1.68%  '
$ws.Range('B33').Value = 'Fetch.AI'
$ws.Range('C33').Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '1.39'
$ws.Range('E33').Value = '  +0.22%  '
$ws.Range('E34').Value = '  +1.46%  '
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D35').Value = '1.00'
$ws.Range('E35').Value = '  -0.07%  '
$ws.Range('E36').Value = '  +2.95%  '
$ws.Range('E37').Value = '  +0.98%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '162.81'
$ws.Range('E38').Value = '  +4.86%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '19.08'
$ws.Range('E39').Value = '  +0.70%  '
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D40').Value = '1.88'
$ws.Range('E40').Value = '  +4.42%  '
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D41').Value = '0.366'
$ws.Range('E41').Value = '  +0.04%  '
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D42').Value = '5.30'
$ws.Range('E42').Value = '  +1.33%  '
$ws.Range('D43').NumberFormat = "@"
$ws.Range('D43').Value = '2.65'
$ws.Range('E43').Value = '  +2.68%  '
$ws.Range('B44').Value = 'WhiteBITCoin'
$ws.Range('C44').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '17.46'
$ws.Range('E44').Value = '  +4.55%  '
$ws.Range('B45').Value = 'BabyDogeCoin'
$ws.Range('C45').Value = 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge'
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '0.0₆0325'
$ws.Range('E45').Value = '  +10.85%  '
$ws.Range('E46').Value = '  +0.14%  '
$ws.Range('D47').NumberFormat = "@"
$ws.Range('D47').Value = '40.10'
$ws.Range('E47').Value = '  -2.01%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '154.69'
$ws.Range('E48').Value = '  +0.07%  '
$ws.Range('E49').Value = '  -1.30%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '21.79'
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('E51').Value = '  -0.71%  '
